# Commit: "Change '_' to '-' in Excel metadata" (Close #5)
#
# The metadata keyword list on the "General_MD" sheet used underscores in
# several of its keys (SUBJECT_CODE, SUBJECT_AREA_en, SUBJECT_AREA_fo,
# CREATION_DATE, UPDATE_FREQUENCY, LAST_UPDATED, NEXT_UPDATE, AXIS_VERSION).
# Those separators are changed to dashes (SUBJECT-CODE, SUBJECT-AREA_en, ...)
# while the language suffix ("_en"/"_fo") keeps its underscore.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General_MD")

# Re-key the eight metadata labels in column A (order chosen to match the
# shared-string append order produced by Excel on save).
$ws.Range("A3").Value  = "AXIS-VERSION"
$ws.Range("A19").Value = "SUBJECT-CODE"
$ws.Range("A24").Value = "CREATION-DATE"
$ws.Range("A25").Value = "UPDATE-FREQUENCY"
$ws.Range("A26").Value = "LAST-UPDATED"
$ws.Range("A27").Value = "NEXT-UPDATE"
$ws.Range("A20").Value = "SUBJECT-AREA_fo"
$ws.Range("A21").Value = "SUBJECT-AREA_en"

# The author ended the edit with the cursor on the General_MD sheet
# (cell A7), which became the active tab; Variables_MD's stored selection
# also moved to K1.
$wb.Worksheets.Item("Variables_MD").Range("K1").Select()

$ws.Activate()
$ws.Range("A7").Select()
